$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openbis-metadata")
$ws.Range("B3").Value = "JJS-MGP90"
